$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Enable (run) test suites E and F: flip Runmode from "N" to "Y"
# for TestCase_F3 (row 4) and TestCase_F4 (row 5).
$ws.Range("C4").Value = "Y"
$ws.Range("C5").Value = "Y"

# Move the active selection to C6, matching the post-edit cursor position.
$ws.Range("C6").Select()
